# Autogenerated on Mon Feb 09 2015 03:30:35 GMT+0000 (Coordinated Universal Time)
#
# Inserts the "MSME definitions" table (Number of employees / Assets /
# Turnover by enterprise-size class) into the Liechtenstein Summary sheet,
# between the existing "Value added to the economy" block (row 21) and the
# trailing source-citation rows (which get pushed down from 26/27 to 32/33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: push everything at/after row 23 down by 6 rows. This carries
# the old rows 26/27 ("SME Performance Review EU" + the long citation) down
# to the new rows 32/33, and leaves a clean blank block at rows 23-28 for
# the new table, matching the target dimension A1:D33.
$ws.Range("A23:A28").EntireRow.Insert()

# Header row of the new table (bold "title" style, like the other header
# rows at B11:D11 / B19:D19).
$ws.Range("B23").Value = "Number of employees"
$ws.Range("B23").Style = "title"
$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("C23").Style = "title"
$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("D23").Style = "title"

# Micro
$ws.Range("A24").Value = "Micro"
$ws.Range("B24").Value = "<10"
$ws.Range("C24").Value = ""
$ws.Range("D24").Value = ""

# Small
$ws.Range("A25").Value = "Small"
$ws.Range("B25").Value = "<50"
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""

# Medium
$ws.Range("A26").Value = "Medium"
$ws.Range("B26").Value = "<250"
$ws.Range("C26").Value = ""
$ws.Range("D26").Value = ""

# Large
$ws.Range("A27").Value = "Large"
$ws.Range("B27").Value = ">249"
$ws.Range("C27").Value = ""
$ws.Range("D27").Value = ""
